# Update "想去人数" (F column) values across the workbook sheets to reflect
# refreshed counts as published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$changes = @{
    6  = 198
    7  = 4573
    8  = 193
    11 = 91
    13 = 695
    14 = 183
    15 = 980
    19 = 68
    20 = 116
    22 = 3505
    23 = 5861
    25 = 32
    29 = 3355
    30 = 360
    32 = 2470
    35 = 125
    36 = 214
    38 = 350
    39 = 125
    40 = 1011
    41 = 908
    42 = 19
    43 = 21
    44 = 47
    45 = 49
    46 = 469
    47 = 63
}
foreach ($row in $changes.Keys) {
    $ws.Cells.Item($row, 6).Value = $changes[$row]
}

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 96

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$changes = @{
    6  = 198
    7  = 4573
    8  = 193
    11 = 96
    12 = 91
    14 = 695
    15 = 183
    16 = 980
    20 = 68
    21 = 116
    23 = 3505
    24 = 5861
    26 = 32
    30 = 3355
    31 = 360
    33 = 2470
    36 = 125
    37 = 214
    39 = 350
    40 = 125
    41 = 1011
    42 = 908
    43 = 19
    44 = 21
    45 = 47
    46 = 49
    47 = 469
    48 = 63
}
foreach ($row in $changes.Keys) {
    $ws.Cells.Item($row, 6).Value = $changes[$row]
}
